{"js": "// Office.js (Word JavaScript API) script.\n// Applies the same edits described by the target OOXML diff:\n//  - Splits several run texts and wraps grammar/spelling \"correction\"\n//    words in <w:proofErr> start/end markers (mirrors what Word inserts\n//    automatically when the proofing pass reflows a paragraph).\n//  - Appends extra sentences to a couple of existing paragraphs.\n//  - Rewrites the \"Discuss planned ... 11 A[bookmark]M ...\" paragraph as\n//    a single contiguous sentence and relocates the _GoBack bookmark.\n//  - Appends a whole new \"14-Jul-2017\" day section at the end of the\n//    report, finishing with a new closing paragraph.\n//\n// Implementation strategy: because every paragraph we touch is made of\n// plain <w:r><w:t>...</w:t></w:r> runs with no run/paragraph formatting,\n// we can reliably reproduce the exact target markup (proofErr tags,\n// preserved spaces, bookmark placement, etc.) by replacing a paragraph's\n// OOXML wholesale (or inserting new OOXML after it) via insertOoxml,\n// rather than trying to recreate the same structure through many small\n// range inserts.\n\nconst NS = 'xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"';\n\n// Wrap one or more <w:p>...</w:p> fragments in the minimal OOXML package\n// envelope that Word.Range/Paragraph.insertOoxml expects.\nfunction pkg(bodyInner) {\n  return (\n    '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n    '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n    '<pkg:xmlData>' +\n    '<w:document ' + NS + '><w:body>' + bodyInner + '</w:body></w:document>' +\n    '</pkg:xmlData></pkg:part></pkg:package>'\n  );\n}\n\nasync function run(context) {\n  const body = context.document.body;\n  const paras = body.paragraphs;\n  paras.load('items');\n  await context.sync();\n\n  // ---- 1. \"After security formalities  allocated seat around 11.45 AM\"\n  const pSecurity = paras.items[5];\n  pSecurity.insertOoxml(\n    pkg(\n      '<w:p><w:r><w:t xml:space=\"preserve\">After security </w:t></w:r>' +\n      '<w:proofErr w:type=\"gramStart\"/>' +\n      '<w:r><w:t>formalities  allocated</w:t></w:r>' +\n      '<w:proofErr w:type=\"gramEnd\"/>' +\n      '<w:r><w:t xml:space=\"preserve\"> seat around 11.45 AM</w:t></w:r></w:p>'\n    ),\n    Word.InsertLocation.replace\n  );\n\n  // ---- 2. \"Copied the copypaste access application ...\"\n  const pCopied = paras.items[10];\n  pCopied.insertOoxml(\n    pkg(\n      '<w:p><w:r><w:t xml:space=\"preserve\">Copied the </w:t></w:r>' +\n      '<w:proofErr w:type=\"spellStart\"/>' +\n      '<w:r><w:t>copypaste</w:t></w:r>' +\n      '<w:proofErr w:type=\"spellEnd\"/>' +\n      '<w:r><w:t xml:space=\"preserve\"> access application and the VBA application into a desktop</w:t></w:r></w:p>'\n    ),\n    Word.InsertLocation.replace\n  );\n\n  // ---- 3. \"The end user ... left around 6.15 PM.\"\n  const p615 = paras.items[12];\n  p615.insertOoxml(\n    pkg(\n      '<w:p><w:r><w:t xml:space=\"preserve\">The end user had to leave. He has to be there for me to continue. </w:t></w:r>' +\n      '<w:proofErr w:type=\"gramStart\"/>' +\n      '<w:r><w:t>So</w:t></w:r>' +\n      '<w:proofErr w:type=\"gramEnd\"/>' +\n      '<w:r><w:t xml:space=\"preserve\"> I couldn\\u2019t continue and left around 6.15 PM.</w:t></w:r></w:p>'\n    ),\n    Word.InsertLocation.replace\n  );\n\n  // ---- 4. \"The end user ... left around 6.45 PM.\"\n  const p645 = paras.items[23];\n  p645.insertOoxml(\n    pkg(\n      '<w:p><w:r><w:t xml:space=\"preserve\">The end user had to leave. He has to be there for me to continue. </w:t></w:r>' +\n      '<w:proofErr w:type=\"gramStart\"/>' +\n      '<w:r><w:t>So</w:t></w:r>' +\n      '<w:proofErr w:type=\"gramEnd\"/>' +\n      '<w:r><w:t xml:space=\"preserve\"> I couldn\\u2019t continue and left around 6.45 PM.</w:t></w:r></w:p>'\n    ),\n    Word.InsertLocation.replace\n  );\n\n  // ---- 5. \"Documented the various testing scenarios ...\" (13-Jul), plus\n  //         the new trailing \" for Maharashtra and Goa circle\".\n  const pScenarios13 = paras.items[29];\n  pScenarios13.insertOoxml(\n    pkg(\n      '<w:p><w:r><w:t>Documented the various testing scenarios</w:t></w:r>' +\n      '<w:r><w:t xml:space=\"preserve\"> \\u2013 Requests without unbilled information, Requests with unbilled information, </w:t></w:r>' +\n      '<w:proofErr w:type=\"gramStart\"/>' +\n      '<w:r><w:t>First</w:t></w:r>' +\n      '<w:proofErr w:type=\"gramEnd\"/>' +\n      '<w:r><w:t xml:space=\"preserve\"> time Credit posting transaction for a customer, Customers with credit posting already done (to check for duplicity)</w:t></w:r>' +\n      '<w:r><w:t xml:space=\"preserve\"> for Maharashtra and Goa circle</w:t></w:r></w:p>'\n    ),\n    Word.InsertLocation.replace\n  );\n\n  // ---- 6. \"Discuss planned with Kunal and Mohit at 11 A[_GoBack]M ...\"\n  //         becomes one contiguous sentence; bookmark removed from here\n  //         (it moves to the new \"Completed the BOLT ...\" paragraph, #8).\n  const pDiscuss = paras.items[36];\n  pDiscuss.insertOoxml(\n    pkg(\n      '<w:p><w:r><w:t>Discuss planned with Kunal and Mohit at 11 AM in ground floor meeting room on the progress and next steps.</w:t></w:r></w:p>'\n    ),\n    Word.InsertLocation.replace\n  );\n\n  // ---- 7. \"The end user ... left around 7 PM.\" gets the same gramStart/\n  //         gramEnd \"So\" split as the other two end-of-day sentences.\n  //         NB: per the source diff this one uses a straight apostrophe.\n  const p7pm = paras.items[37];\n  p7pm.insertOoxml(\n    pkg(\n      '<w:p><w:r><w:t xml:space=\"preserve\">The end user had to leave. He has to be there for me to continue. </w:t></w:r>' +\n      '<w:proofErr w:type=\"gramStart\"/>' +\n      '<w:r><w:t>So</w:t></w:r>' +\n      '<w:proofErr w:type=\"gramEnd\"/>' +\n      \"<w:r><w:t xml:space=\\\"preserve\\\"> I couldn't continue and left around 7 PM.</w:t></w:r></w:p>\"\n    ),\n    Word.InsertLocation.replace\n  );\n\n  // ---- 8. Append the whole new \"14-Jul-2017\" section after paragraph 37.\n  const newDaySection =\n    '<w:p><w:r><w:br w:type=\"page\"/></w:r></w:p>' +\n    '<w:p><w:r><w:lastRenderedPageBreak/><w:t>14-Jul-2017</w:t></w:r></w:p>' +\n    '<w:p><w:r><w:t xml:space=\"preserve\">Arrived in Vodafone office at </w:t></w:r>' +\n    '<w:r><w:t>9.45</w:t></w:r>' +\n    '<w:r><w:t xml:space=\"preserve\"> AM</w:t></w:r></w:p>' +\n    '<w:p><w:r><w:t>Work</w:t></w:r>' +\n    '<w:r><w:t xml:space=\"preserve\"> delayed by 30 minutes </w:t></w:r>' +\n    '<w:r><w:t>for waiting for the Vodafone employee to unlock the system</w:t></w:r></w:p>' +\n    '<w:p><w:r><w:t xml:space=\"preserve\">Documented the various testing scenarios \\u2013 Requests without unbilled information, Requests with unbilled information, </w:t></w:r>' +\n    '<w:proofErr w:type=\"gramStart\"/>' +\n    '<w:r><w:t>First</w:t></w:r>' +\n    '<w:proofErr w:type=\"gramEnd\"/>' +\n    '<w:r><w:t xml:space=\"preserve\"> time Credit posting transaction for a customer, Customers with credit posting already done (to check for duplicity)</w:t></w:r>' +\n    '<w:r><w:t xml:space=\"preserve\"> for a different circle (Mumbai)</w:t></w:r></w:p>' +\n    '<w:p><w:r><w:t xml:space=\"preserve\">Updated the screenshot document by formatting the document, giving captions to screen shots  </w:t></w:r></w:p>' +\n    '<w:p><w:r><w:t xml:space=\"preserve\">Worked on the </w:t></w:r>' +\n    '<w:r><w:t>Tool Scenarios, Results, Process Logic and Observations</w:t></w:r>' +\n    '<w:r><w:t xml:space=\"preserve\">, </w:t></w:r>' +\n    '<w:r><w:t>Test conditions handled by the tool</w:t></w:r>' +\n    '<w:r><w:t xml:space=\"preserve\"> documents</w:t></w:r></w:p>' +\n    '<w:p><w:r><w:t>Completed</w:t></w:r>' +\n    '<w:r><w:t xml:space=\"preserve\"> the </w:t></w:r>' +\n    '<w:r><w:t>BOLT Application review document</w:t></w:r>' +\n    '<w:r><w:t xml:space=\"preserve\"> \\u2013 the template provided by Vodafone to fill the observations at appropriate </w:t></w:r>' +\n    '<w:r><w:t>places thereby completing the assessment process.</w:t></w:r>' +\n    '<w:bookmarkStart w:id=\"0\" w:name=\"_GoBack\"/><w:bookmarkEnd w:id=\"0\"/></w:p>';\n\n  const rangeAfter7pm = p7pm.getRange();\n  rangeAfter7pm.insertOoxml(pkg(newDaySection), Word.InsertLocation.after);\n\n  await context.sync();\n\n  // ---- 9. Replace the trailing empty paragraph with the new closing\n  //         sentence (the last paragraph mark of the body can't be\n  //         deleted outright, so we overwrite its contents instead).\n  paras.load('items');\n  await context.sync();\n  const pLast = paras.items[paras.items.length - 1];\n  pLast.insertOoxml(\n    pkg(\n      '<w:p><w:r><w:t>Had a formal meeting with Mohit and Kunal on the progress and suggested action steps. This will be discussed by the management for the action plan next week.</w:t></w:r></w:p>'\n    ),\n    Word.InsertLocation.replace\n  );\n\n  await context.sync();\n}\n\nawait run(context);\n", "ps1": "# Word COM interop (PowerShell-style) script.\n# Applies the same edits described by the target OOXML diff:\n#  - Splits several run texts and wraps grammar/spelling \"correction\"\n#    words in <w:proofErr> start/end markers (mirrors what Word inserts\n#    automatically when the proofing pass reflows a paragraph).\n#  - Appends extra sentences to a couple of existing paragraphs.\n#  - Rewrites the \"Discuss planned ... 11 A[bookmark]M ...\" paragraph as\n#    a single contiguous sentence and relocates the _GoBack bookmark.\n#  - Appends a whole new \"14-Jul-2017\" day section at the end of the\n#    report, finishing with a new closing paragraph.\n#\n# Implementation strategy: every paragraph we touch is made of plain\n# runs with no run/paragraph formatting, so each edit is applied by\n# feeding a Range.InsertXML() call the exact WordprocessingML we want in\n# place (wrapped in the minimal xmlPackage envelope Word expects). Using\n# a paragraph's full Range performs a clean \"replace this paragraph's\n# content\" - the same effect as Office.js's InsertLocation.Replace.\n# New paragraphs are added by collapsing the Range of the document's\n# final (always-empty) paragraph to its Start and inserting there, which\n# pushes the new paragraphs in front of it without disturbing anything\n# already in the document or losing that trailing paragraph mark.\n\n$d = $word.ActiveDocument\n\nfunction New-PkgXml([string]$bodyInner) {\n    return '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n        '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n        '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n        '<pkg:xmlData>' +\n        '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"><w:body>' + $bodyInner + '</w:body></w:document>' +\n        '</pkg:xmlData></pkg:part></pkg:package>'\n}\n\n# ---- 1. \"After security formalities  allocated seat around 11.45 AM\" (para 6)\n$r = $d.Paragraphs.Item(6).Range\n$r.InsertXML((New-PkgXml (\n    '<w:p><w:r><w:t xml:space=\"preserve\">After security </w:t></w:r>' +\n    '<w:proofErr w:type=\"gramStart\"/>' +\n    '<w:r><w:t>formalities  allocated</w:t></w:r>' +\n    '<w:proofErr w:type=\"gramEnd\"/>' +\n    '<w:r><w:t xml:space=\"preserve\"> seat around 11.45 AM</w:t></w:r></w:p>'\n)))\n\n# ---- 2. \"Copied the copypaste access application ...\" (para 11)\n$r = $d.Paragraphs.Item(11).Range\n$r.InsertXML((New-PkgXml (\n    '<w:p><w:r><w:t xml:space=\"preserve\">Copied the </w:t></w:r>' +\n    '<w:proofErr w:type=\"spellStart\"/>' +\n    '<w:r><w:t>copypaste</w:t></w:r>' +\n    '<w:proofErr w:type=\"spellEnd\"/>' +\n    '<w:r><w:t xml:space=\"preserve\"> access application and the VBA application into a desktop</w:t></w:r></w:p>'\n)))\n\n# ---- 3. \"The end user ... left around 6.15 PM.\" (para 13)\n$r = $d.Paragraphs.Item(13).Range\n$r.InsertXML((New-PkgXml (\n    '<w:p><w:r><w:t xml:space=\"preserve\">The end user had to leave. He has to be there for me to continue. </w:t></w:r>' +\n    '<w:proofErr w:type=\"gramStart\"/>' +\n    '<w:r><w:t>So</w:t></w:r>' +\n    '<w:proofErr w:type=\"gramEnd\"/>' +\n    '<w:r><w:t xml:space=\"preserve\"> I couldn' + [char]0x2019 + 't continue and left around 6.15 PM.</w:t></w:r></w:p>'\n)))\n\n# ---- 4. \"The end user ... left around 6.45 PM.\" (para 24)\n$r = $d.Paragraphs.Item(24).Range\n$r.InsertXML((New-PkgXml (\n    '<w:p><w:r><w:t xml:space=\"preserve\">The end user had to leave. He has to be there for me to continue. </w:t></w:r>' +\n    '<w:proofErr w:type=\"gramStart\"/>' +\n    '<w:r><w:t>So</w:t></w:r>' +\n    '<w:proofErr w:type=\"gramEnd\"/>' +\n    '<w:r><w:t xml:space=\"preserve\"> I couldn' + [char]0x2019 + 't continue and left around 6.45 PM.</w:t></w:r></w:p>'\n)))\n\n# ---- 5. \"Documented the various testing scenarios ...\" (13-Jul, para 30),\n#         plus the new trailing \" for Maharashtra and Goa circle\".\n$r = $d.Paragraphs.Item(30).Range\n$r.InsertXML((New-PkgXml (\n    '<w:p><w:r><w:t>Documented the various testing scenarios</w:t></w:r>' +\n    '<w:r><w:t xml:space=\"preserve\"> ' + [char]0x2013 + ' Requests without unbilled information, Requests with unbilled information, </w:t></w:r>' +\n    '<w:proofErr w:type=\"gramStart\"/>' +\n    '<w:r><w:t>First</w:t></w:r>' +\n    '<w:proofErr w:type=\"gramEnd\"/>' +\n    '<w:r><w:t xml:space=\"preserve\"> time Credit posting transaction for a customer, Customers with credit posting already done (to check for duplicity)</w:t></w:r>' +\n    '<w:r><w:t xml:space=\"preserve\"> for Maharashtra and Goa circle</w:t></w:r></w:p>'\n)))\n\n# ---- 6. \"Discuss planned with Kunal and Mohit at 11 A[_GoBack]M ...\"\n#         becomes one contiguous sentence; bookmark removed from here\n#         (it moves to the new \"Completed the BOLT ...\" paragraph, #8).\n$r = $d.Paragraphs.Item(37).Range\n$r.InsertXML((New-PkgXml (\n    '<w:p><w:r><w:t>Discuss planned with Kunal and Mohit at 11 AM in ground floor meeting room on the progress and next steps.</w:t></w:r></w:p>'\n)))\n\n# ---- 7. \"The end user ... left around 7 PM.\" gets the same gramStart/\n#         gramEnd \"So\" split as the other two end-of-day sentences.\n#         NB: per the source diff this one uses a straight apostrophe.\n$r = $d.Paragraphs.Item(38).Range\n$r.InsertXML((New-PkgXml (\n    '<w:p><w:r><w:t xml:space=\"preserve\">The end user had to leave. He has to be there for me to continue. </w:t></w:r>' +\n    '<w:proofErr w:type=\"gramStart\"/>' +\n    '<w:r><w:t>So</w:t></w:r>' +\n    '<w:proofErr w:type=\"gramEnd\"/>' +\n    '<w:r><w:t xml:space=\"preserve\"> I couldn' + [char]0x27 + 't continue and left around 7 PM.</w:t></w:r></w:p>'\n)))\n\n# ---- 8. Append the whole new \"14-Jul-2017\" section, and ---- 9. the new\n#         closing sentence, by inserting before the document's always-\n#         empty trailing paragraph (its Start never moves underneath us).\n$newDaySection =\n    '<w:p><w:r><w:br w:type=\"page\"/></w:r></w:p>' +\n    '<w:p><w:r><w:lastRenderedPageBreak/><w:t>14-Jul-2017</w:t></w:r></w:p>' +\n    '<w:p><w:r><w:t xml:space=\"preserve\">Arrived in Vodafone office at </w:t></w:r>' +\n    '<w:r><w:t>9.45</w:t></w:r>' +\n    '<w:r><w:t xml:space=\"preserve\"> AM</w:t></w:r></w:p>' +\n    '<w:p><w:r><w:t>Work</w:t></w:r>' +\n    '<w:r><w:t xml:space=\"preserve\"> delayed by 30 minutes </w:t></w:r>' +\n    '<w:r><w:t>for waiting for the Vodafone employee to unlock the system</w:t></w:r></w:p>' +\n    '<w:p><w:r><w:t xml:space=\"preserve\">Documented the various testing scenarios ' + [char]0x2013 + ' Requests without unbilled information, Requests with unbilled information, </w:t></w:r>' +\n    '<w:proofErr w:type=\"gramStart\"/>' +\n    '<w:r><w:t>First</w:t></w:r>' +\n    '<w:proofErr w:type=\"gramEnd\"/>' +\n    '<w:r><w:t xml:space=\"preserve\"> time Credit posting transaction for a customer, Customers with credit posting already done (to check for duplicity)</w:t></w:r>' +\n    '<w:r><w:t xml:space=\"preserve\"> for a different circle (Mumbai)</w:t></w:r></w:p>' +\n    '<w:p><w:r><w:t xml:space=\"preserve\">Updated the screenshot document by formatting the document, giving captions to screen shots  </w:t></w:r></w:p>' +\n    '<w:p><w:r><w:t xml:space=\"preserve\">Worked on the </w:t></w:r>' +\n    '<w:r><w:t>Tool Scenarios, Results, Process Logic and Observations</w:t></w:r>' +\n    '<w:r><w:t xml:space=\"preserve\">, </w:t></w:r>' +\n    '<w:r><w:t>Test conditions handled by the tool</w:t></w:r>' +\n    '<w:r><w:t xml:space=\"preserve\"> documents</w:t></w:r></w:p>' +\n    '<w:p><w:r><w:t>Completed</w:t></w:r>' +\n    '<w:r><w:t xml:space=\"preserve\"> the </w:t></w:r>' +\n    '<w:r><w:t>BOLT Application review document</w:t></w:r>' +\n    '<w:r><w:t xml:space=\"preserve\"> ' + [char]0x2013 + ' the template provided by Vodafone to fill the observations at appropriate </w:t></w:r>' +\n    '<w:r><w:t>places thereby completing the assessment process.</w:t></w:r>' +\n    '<w:bookmarkStart w:id=\"0\" w:name=\"_GoBack\"/><w:bookmarkEnd w:id=\"0\"/></w:p>'\n\n$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)\n$rEnd = $lastPara.Range\n$rEnd.Collapse(1)\n$rEnd.InsertXML((New-PkgXml $newDaySection))\n\n# ---- 9. Fill in the (still) trailing empty paragraph with the new\n#         closing sentence. A plain Range.Text assignment (rather than\n#         InsertXML) is used here because this final paragraph mark is\n#         the document body's terminating one and can never truly be\n#         replaced wholesale - InsertXML-ing a full <w:p> \"over\" it\n#         leaves a stray empty paragraph behind, whereas setting .Text\n#         reuses the existing mark cleanly (the text itself needs no\n#         run splitting/proofErr markers, so plain text is sufficient).\n$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)\n$rEnd = $lastPara.Range\n$rEnd.Text = \"Had a formal meeting with Mohit and Kunal on the progress and suggested action steps. This will be discussed by the management for the action plan next week.\"\n"}
